$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tasks appended below the existing list (kept on odd rows, matching
# the sheet's existing one-blank-row-between-items layout).
$ws.Range("B27").Value = "Read about  Owin middleware"
$ws.Range("B29").Value = "Implement membership Service, CryptoService and everything for authentication and users"
$ws.Range("B31").Value = "Complete membership and User Services and Managers "

# Select the newly added last item, and scroll the viewport down so row 17
# is at the top-left of the visible pane (mirrors the author scrolling down
# to see their new entries after typing them).
$ws.Range("B31").Select()
try {
    $excel.ActiveWindow.ScrollRow = 17
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Older/headless hosts may not support window-scroll properties; the
    # selection above is the functionally important part.
}
